# Automatische test-sync: 2025-08-03 18:42:50
# Adds Testmail #14 ("Heb je de CE-certificaten van dit product?") to the
# Logs sheet (row 42) and the matching tally row to the Dashboard sheet
# (row 10), then extends the conditional formatting ranges and the chart's
# category/value series references to cover the new row.

$wb = $excel.ActiveWorkbook

# --- 1. Logs sheet: append row 42 -----------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(42, 1).Value  = "Heb je de CE-certificaten van dit product?"
$logs.Cells.Item(42, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item(42, 3).Value  = "Testmail #14: Heb je de CE-certificaten van dit product?"
$logs.Cells.Item(42, 4).Value  = "Kwaliteit / Certificaten"
$logs.Cells.Item(42, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar kwaliteit@bedrijf.nl."
$logs.Cells.Item(42, 6).Value  = "2025-08-03 18:42:24"
$logs.Cells.Item(42, 7).Value  = "Ja"
$logs.Cells.Item(42, 8).Value  = "Ja"
$logs.Cells.Item(42, 9).Value  = "Nee"
$logs.Cells.Item(42, 10).Value = "Nee"

# --- 2. Dashboard sheet: append tally row 10 -------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(10, 1).Value = "Kwaliteit / Certificaten"
$dashboard.Cells.Item(10, 2).Value = 1

# --- 3. Grow the conditional formatting ranges on Logs by one row ---------
$logs.Range("D2:D41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D42"))
$logs.Range("G2:G41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G42"))
$logs.Range("H2:H41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H42"))
$logs.Range("I2:I41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I42"))
$logs.Range("J2:J41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J42"))

# --- 4. Extend the Dashboard bar chart's category/value series ------------
$chartObj = $dashboard.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$10,Dashboard!`$B`$2:`$B`$10,1)"
